# "shift to prod env"
#
# The ScenarioMapping sheet drives D336:D370 (the "SmokeTest" column) from
# "No" to "Yes" for that block of scenarios, and leaves the window scrolled
# near the bottom of the sheet with F369 as the active cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the SmokeTest flag (column D) from "No" to "Yes" for rows 336-370.
$rng = $ws.Range("D336:D370")
for ($i = 1; $i -le $rng.Rows.Count; $i++) {
    $rng.Cells.Item($i, 1).Value = "Yes"
}

# Leave the sheet scrolled down near the bottom with F369 selected, as the
# view was left after making the edits.
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 360
$ws.Range("F369").Select()
